$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 2 ("Ecriture de la documentation"): duration grows from 2 to 3 days
# -> its END DATE formula (E2 = C2+D2) recalculates automatically.
$ws.Range("D2").Value = 3

# Row 4 ("Creation de la maquette du site web"): the task now starts a day
# later and gets a 2-day duration -> END DATE formula (E4) recalculates.
$ws.Range("C4").Value = 42123
$ws.Range("D4").Value = 2

# Bring the view back to the top of the sheet and move the active selection
# to I9 (previously the view was scrolled down with G9 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I9").Select()

$wb.Save()
